$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.681.96'
$ws.Range("E2").Value = '  -4.14%  '
$ws.Range("D3").Value = '3.118.70'
$ws.Range("E3").Value = '  -3.86%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''553.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.39%  '
$ws.Range("D6").Value = '''138.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -11.10%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.113.11'
$ws.Range("E8").Value = '  -3.77%  '
$ws.Range("D9").Value = '''0.500'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.19%  '
$ws.Range("D10").Value = '''0.158'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.88%  '
$ws.Range("D11").Value = '''6.41'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -9.06%  '
$ws.Range("D12").Value = '''0.474'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.13%  '
$ws.Range("D13").Value = '''35.80'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.66%  '
$ws.Range("D14").Value = '''0.0000220'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.09%  '
$ws.Range("D15").Value = '3.598.12'
$ws.Range("E15").Value = '  -4.34%  '
$ws.Range("D16").Value = '63.551.21'
$ws.Range("E16").Value = '  -4.45%  '
$ws.Range("D17").Value = '''0.112'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.15%  '
$ws.Range("D18").Value = '3.094.81'
$ws.Range("E18").Value = '  -4.41%  '
$ws.Range("D19").Value = '''6.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.76%  '
$ws.Range("D20").Value = '''495.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -11.92%  '
$ws.Range("D21").Value = '''13.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.73%  '
$ws.Range("D22").Value = '''0.727'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.59%  '
$ws.Range("D23").Value = '''7.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.23%  '
$ws.Range("D24").Value = '''79.61'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.37%  '
$ws.Range("D25").Value = '''12.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.48%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").Value = '''8.53'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -9.72%  '
$ws.Range("D28").Value = '''2.77'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.09%  '
$ws.Range("D29").Value = '''2.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -11.96%  '
$ws.Range("D30").Value = '''0.996'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.45%  '
$ws.Range("D31").Value = '''26.82'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.27%  '
$ws.Range("D32").Value = '''1.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.43%  '
$ws.Range("D33").Value = '''2.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.67%  '
$ws.Range("D34").Value = '''59.15'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.70%  '
$ws.Range("D35").Value = '''521.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.77%  '
$ws.Range("D36").Value = '''6.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.73%  '
$ws.Range("D37").Value = '''5.21'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -10.15%  '
$ws.Range("D38").Value = '''0.0410'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -11.40%  '
$ws.Range("D39").Value = '3.160.99'
$ws.Range("E39").Value = '  +0.56%  '
$ws.Range("D40").Value = '''0.0810'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.92%  '
$ws.Range("E41").Value = '  -5.65%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '''2.71'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -11.46%  '
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").Value = '''8.22'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.24%  '
$ws.Range("D44").Value = '''0.261'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.43%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").Value = '''2.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.56%  '
$ws.Range("D47").Value = '''25.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.08%  '
$ws.Range("D48").Value = '''121.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("E49").Value = '  -3.90%  '
$ws.Range("D50").Value = '0.0₃0512'
$ws.Range("E50").Value = '  -9.22%  '
$ws.Range("D51").Value = '''2.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.55%  '
